# updated main GSC export data
#
# Append the next day's row (2025-12-12) to the "Chart" data sheet, carrying
# forward the same HTTPS URL count (29) and a flat 0 for Non-HTTPS URLs -
# mirroring the previous day's row (2025-12-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Force the date column to be stored as literal text (matching every other
# row in the export, which are plain text dates, not Excel date serials).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-12"

# Copying the previous row's formatting back over the new cell collapses it
# onto the workbook's existing (default) style instead of leaving behind the
# one-off "@" text format, just like the rest of the sheet.
$ws.Cells.Item($lastRow, 1).Copy()
$dateCell.PasteSpecial(-4122)

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 29
